$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.777.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "'1.649.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'215.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +2.32%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'19.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "'1.635.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").Value = "'0.532"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "'66.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.20%  "
$ws.Range("D17").Value = "'26.817.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'0.0₃0747"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "'218.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.37%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").Value = "'6.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.17%  "
$ws.Range("D23").Value = "'9.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.60%  "
$ws.Range("D25").Value = "'147.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").Value = "'6.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "'15.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("D33").Value = "'3.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("D34").Value = "'1.279.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.91%  "
$ws.Range("D35").Value = "'1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.0179"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.07%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.517"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.810"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'2.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").Value = "'5.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "'1.787.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "'93.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("E46").Value = "  +4.34%  "
$ws.Range("D47").Value = "'56.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("D51").Value = "'0.0969"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.19%  "
